# Revert "Powerpoint writer: consolidate text run nodes."
# Re-split title runs that were merged back into "word + trailing space"
# single runs, so each word and the following space become separate
# <a:r> runs again (matching the original, pre-consolidation output).

$p = $ppt.ActivePresentation

# Slide 1 title: "Header " / "with " -> "Header" / " " / "with" / " "
$s1 = $p.Slides.Item(1)
$tr1 = $s1.Shapes.Item(1).TextFrame.TextRange
$tr1.Characters(1, 6).Text = "Header"
$tr1.Characters(8, 4).Text = "with"

# Slide 2 title: "Syntax " -> "Syntax" / " "
$s2 = $p.Slides.Item(2)
$tr2 = $s2.Shapes.Item(1).TextFrame.TextRange
$tr2.Characters(1, 6).Text = "Syntax"

# Slide 3 title: "Two " / "column " -> "Two" / " " / "column" / " "
$s3 = $p.Slides.Item(3)
$tr3 = $s3.Shapes.Item(1).TextFrame.TextRange
$tr3.Characters(1, 3).Text = "Two"
$tr3.Characters(5, 6).Text = "column"
